$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add notes for the newly completed tasks (order matters for shared-string indices)
$ws.Range("G35").Value = "Camera shifts down when launched up on jump pad"
$ws.Range("G36").Value = "Camera shakes slightly when player is moving fast"
$ws.Range("G30").Value = "Made jump pad"
$ws.Range("G20").Value = "Changed textures on all walls, changed pickup look, added red lighting"

# Mark the corresponding checkboxes as completed (linked boolean cells)
$ws.Range("J20").Value = $true
$ws.Range("J30").Value = $true
$ws.Range("J35").Value = $true
$ws.Range("J36").Value = $true

# Update the last active selection to reflect where the user left off
$ws.Range("G20").Select()
